$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 40, shifting existing rows 40-108 down to 42-110
$ws.Rows("40:41").Insert()

# Populate the two new rows (40 and 41) with the new data record
# Row 40
$ws.Range("A40").Value = 2
$ws.Range("B40").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C40").Value = 'Coquimbo'
$ws.Range("D40").Value = 44994
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = 'Fruta'
$ws.Range("G40").Value = 100103
$ws.Range("H40").Value = 'Frutos de hueso (carozo)'
$ws.Range("I40").Value = 100103002
$ws.Range("J40").Value = 'Ciruela'
$ws.Range("K40").Value = 'Angeleno'
$ws.Range("L40").Value = 'Primera'
$ws.Range("M40").Value = 10
$ws.Range("N40").Value = 200000
$ws.Range("O40").Value = 210000
$ws.Range("P40").Value = 205000
$ws.Range("Q40").Value = '$/bins (450 kilos)'
$ws.Range("R40").Value = 'Región de O''Higgins'
$ws.Range("S40").Value = 456
$ws.Range("T40").Value = 450

# Row 41
$ws.Range("A41").Value = 2
$ws.Range("B41").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C41").Value = 'Coquimbo'
$ws.Range("D41").Value = 44994
$ws.Range("E41").Value = 4
$ws.Range("F41").Value = 'Fruta'
$ws.Range("G41").Value = 100103
$ws.Range("H41").Value = 'Frutos de hueso (carozo)'
$ws.Range("I41").Value = 100103002
$ws.Range("J41").Value = 'Ciruela'
$ws.Range("K41").Value = 'Angeleno'
$ws.Range("L41").Value = 'Segunda'
$ws.Range("M41").Value = 16
$ws.Range("N41").Value = 170000
$ws.Range("O41").Value = 180000
$ws.Range("P41").Value = 175000
$ws.Range("Q41").Value = '$/bins (450 kilos)'
$ws.Range("R41").Value = 'Región de O''Higgins'
$ws.Range("S41").Value = 389
$ws.Range("T41").Value = 450
